# Update with Correct Forecast output
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before column B; this shifts ASIN..is_holiday_week one
# column to the right (B->C, C->D, D->E, E->F, F->G, G->H, H->I, I->J).
$ws.Columns.Item(2).Insert()

# Header for the newly-inserted column
$ws.Range("B1").Value = "Week_Start_Date"

# Week start dates for rows 2..17 (weekly cadence starting 2025-01-05),
# stored as literal text (not auto-converted Excel date serials).
$dates = @(
    "2025-01-05", "2025-01-12", "2025-01-19", "2025-01-26",
    "2025-02-02", "2025-02-09", "2025-02-16", "2025-02-23",
    "2025-03-02", "2025-03-09", "2025-03-16", "2025-03-23",
    "2025-03-30", "2025-04-06", "2025-04-13", "2025-04-20"
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Range("B" + $row)
    $cell.NumberFormat = "@"
    $cell.Value = $dates[$i]
    $cell.Style = "Normal"
}

# Fix Week labels in column A: strip leading zero, W01..W09 -> W1..W9
for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Range("A" + $row)
    $val = $cell.Value2
    if ($val -match "^W0(\d)$") {
        $cell.Value = "W" + $matches[1]
    }
}

# is_holiday_week column (now column J) becomes boolean FALSE instead of numeric 0
for ($row = 2; $row -le 17; $row++) {
    $ws.Range("J" + $row).Value = $false
}
